$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 43

# Column A holds a date-like string ("2025-09-27") that must be stored as
# literal text, not auto-converted into a date serial number. Temporarily
# force a text number format while assigning the value, then restore the
# cell to the unstyled "Normal" style so no extra formatting is left behind
# (matching the other data rows, which carry no style attribute).
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-09-27"
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 55.27999877929688
$ws.Cells.Item($row, 3).Value = 672.9000244140625
$ws.Cells.Item($row, 4).Value = 321
